# i18n(CWL): update loc entries
# Adds a new localization row (id + JP/EN text) to the "General" sheet of
# cwl_sources.xlsx:
#   cwl_log_post_cleanup_quest  ->  "removed invalid quest id: {0}"
# and updates the sheet view's scroll/selection position to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 50) so the new row matches the
# existing look (font, wrap, vertical alignment) instead of the empty
# placeholder style that row 52 currently has.
$ws.Range("A50").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("C50").Copy()
$ws.Range("C52").PasteSpecial(-4122)
$ws.Range("D50").Copy()
$ws.Range("D52").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New id (column A) and JP/EN message (columns C & D) for row 52
$ws.Range("A52").Value = "cwl_log_post_cleanup_quest"
$ws.Range("C52").Value = "removed invalid quest id: {0}"
$ws.Range("D52").Value = "removed invalid quest id: {0}"

# Update the view so the newly added row is the active selection
$ws.Activate()
$ws.Range("A52").Select()
